$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''63.646.44'
$ws.Range("E2").Value = '  -0.48%  '
$ws.Range("D3").Value = '''3.090.22'
$ws.Range("E3").Value = '  -1.62%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").Value = '''592.00'
$ws.Range("E5").Value = '  +0.28%  '
$ws.Range("D6").Value = '''155.20'
$ws.Range("E6").Value = '  +6.60%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = '''0.545'
$ws.Range("E8").Value = '  +2.97%  '
$ws.Range("D9").Value = '''3.083.68'
$ws.Range("E9").Value = '  -1.57%  '
$ws.Range("D10").Value = '''0.157'
$ws.Range("E10").Value = '  -1.55%  '
$ws.Range("E11").Value = '  -0.90%  '
$ws.Range("D12").Value = '''0.457'
$ws.Range("E12").Value = '  -0.39%  '
$ws.Range("D13").Value = '''37.59'
$ws.Range("E13").Value = '  +0.63%  '
$ws.Range("E14").Value = '  -1.85%  '
$ws.Range("D15").Value = '''3.601.79'
$ws.Range("E15").Value = '  -1.62%  '
$ws.Range("E16").Value = '  -1.50%  '
$ws.Range("E17").Value = '  -2.17%  '
$ws.Range("D18").Value = '''63.600.46'
$ws.Range("E18").Value = '  -0.31%  '
$ws.Range("D19").Value = '''3.088.48'
$ws.Range("E19").Value = '  -1.60%  '
$ws.Range("D20").Value = '''476.87'
$ws.Range("E20").Value = '  +1.83%  '
$ws.Range("E21").Value = '  +2.22%  '
$ws.Range("D22").Value = '''0.717'
$ws.Range("E22").Value = '  -2.04%  '
$ws.Range("D23").Value = '''7.56'
$ws.Range("E23").Value = '  +0.21%  '
$ws.Range("D24").Value = '''2.41'
$ws.Range("E24").Value = '  +4.20%  '
$ws.Range("D25").Value = '''12.95'
$ws.Range("E25").Value = '  -0.25%  '
$ws.Range("D26").Value = '''81.29'
$ws.Range("E26").Value = '  -0.43%  '
$ws.Range("D27").Value = '''10.06'
$ws.Range("E27").Value = '  +2.52%  '
$ws.Range("E29").Value = '  -0.72%  '
$ws.Range("E30").Value = '  -0.55%  '
$ws.Range("E31").Value = '  +0.01%  '
$ws.Range("D32").Value = '''2.18'
$ws.Range("E32").Value = '  -2.49%  '
$ws.Range("E33").Value = '  +4.05%  '
$ws.Range("E34").Value = '  -1.54%  '
$ws.Range("D35").Value = '''0.0₃0853'
$ws.Range("E35").Value = '  +1.00%  '
$ws.Range("E36").Value = '  -1.29%  '
$ws.Range("D37").Value = '''3.41'
$ws.Range("E37").Value = '  +6.32%  '
$ws.Range("D38").Value = '''6.10'
$ws.Range("E38").Value = '  -0.98%  '
$ws.Range("E39").Value = '  -3.05%  '
$ws.Range("D40").Value = '''9.34'
$ws.Range("E40").Value = '  -0.46%  '
$ws.Range("D41").Value = '''50.70'
$ws.Range("E41").Value = '  -1.35%  '
$ws.Range("D42").Value = '''444.34'
$ws.Range("E42").Value = '  -2.16%  '
$ws.Range("E43").Value = '  -2.19%  '
$ws.Range("D44").Value = '''0.0364'
$ws.Range("E44").Value = '  -2.20%  '
$ws.Range("D45").Value = '''40.03'
$ws.Range("E45").Value = '  +0.95%  '
$ws.Range("E46").Value = '  +2.97%  '
$ws.Range("D47").Value = '''2.804.35'
$ws.Range("E47").Value = '  -3.88%  '
$ws.Range("D48").Value = '''131.75'
$ws.Range("E48").Value = '  -1.05%  '
$ws.Range("D49").Value = '''25.75'
$ws.Range("E49").Value = '  +6.58%  '
$ws.Range("E50").Value = '  +0.04%  '
$ws.Range("E51").Value = '  +1.15%  '
